# Trade #28 closed at 2026-02-16 21:26:55 - leadlag DOWN +0.000%
# Appends a new "OPEN" trade row (row 24) to the "leadlag" sheet,
# mirroring the existing trade-log rows (A1:N23 -> A1:N24).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 24

$ws.Cells.Item($row, 1).Value = 28

# Column B holds a date-like label ("2026-02-16") that must stay plain
# text, not get auto-converted into a date serial number. Force the
# cell to Text first, write the value, then drop the explicit number
# format again so the saved cell carries no extra style (matches the
# rest of the sheet, which relies on the default/general style).
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2026-02-16"
$ws.Cells.Item($row, 2).ClearFormats()

$ws.Cells.Item($row, 3).Value = "21:26:55"
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "DOWN"
$ws.Cells.Item($row, 6).Value = 68907.645
$ws.Cells.Item($row, 7).Value = ""
$ws.Cells.Item($row, 8).Value = "OPEN"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0.75
$ws.Cells.Item($row, 12).Value = "Coinbase leading with -0.169% move"
$ws.Cells.Item($row, 13).Value = ""
$ws.Cells.Item($row, 14).Value = 0
